$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 840, shifting rows 840:881 down to 841:882.
$ws.Rows.Item(840).Insert()

# Populate the newly inserted row with the new daily entry.
# Column A holds a date-like string ("2026/02/19") that must be stored
# verbatim as text (matching the other rows), not auto-converted into a
# date serial number. Temporarily force text format while assigning it,
# then restore the default "Normal" style so no extra formatting lingers
# on the cell.
$ws.Cells.Item(840, 1).NumberFormat = "@"
$ws.Cells.Item(840, 1).Value = "2026/02/19"
$ws.Cells.Item(840, 1).Style = "Normal"
$ws.Cells.Item(840, 2).Value = "木"
$ws.Cells.Item(840, 3).Value = 10
$ws.Cells.Item(840, 4).Value = 201
